{"js": "// 1. Update the \"saved\" date text in the header field to August 18, 2021\nconst dateResults = context.document.body.search(\"February 11, 2021\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"August 18, 2021\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. Insert a new \"South Korea\" hotline row right before the \"Taiwan\" row\n// (the APAC hotline table is the first table in the document)\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst apacTable = tables.items[0];\nconst rows = apacTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\nrows.items.forEach((r) => r.load(\"values\"));\nawait context.sync();\n\nlet taiwanRow = null;\nfor (const r of rows.items) {\n  if (r.values[0][0] === \"Taiwan\") {\n    taiwanRow = r;\n    break;\n  }\n}\n\nif (taiwanRow) {\n  const taiwanValues = taiwanRow.values[0]; // [\"Taiwan\", \"(+886) 00801 852 057\"]\n\n  // Insert the new row *after* Taiwan so it inherits Taiwan's (non-highlighted)\n  // formatting instead of the preceding (highlighted) Singapore row, carrying\n  // Taiwan's own text into that new row...\n  taiwanRow.insertRows(Word.InsertLocation.after, 1, [taiwanValues]);\n  await context.sync();\n\n  // ...then overwrite the original Taiwan row (now duplicated) with the South\n  // Korea entry, so the net effect is \"South Korea\" ends up positioned right\n  // before \"Taiwan\".\n  const taiwanCells = taiwanRow.cells;\n  taiwanCells.load(\"items\");\n  await context.sync();\n\n  taiwanCells.items[0].value = \"South Korea\";\n  taiwanCells.items[1].value = \"(+82) 798 611 4890\";\n  await context.sync();\n}\n\n// 3. Simplify the \"Vietnam (Mobifone)\" country cell to a single run\n// (collapses the proofErr-wrapped, multi-run \"Vietnam (\" + \"Mobifone\" + \")\" text)\nconst vnResults = context.document.body.search(\"Vietnam (Mobifone)\", { matchCase: true });\nvnResults.load(\"items\");\nawait context.sync();\n\nif (vnResults.items.length > 0) {\n  const vnParagraphs = vnResults.items[0].paragraphs;\n  vnParagraphs.load(\"items\");\n  await context.sync();\n\n  const vnRange = vnParagraphs.items[0].getRange();\n  vnRange.insertText(\"Vietnam (Mobifone)\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Update the \"saved\" date text in the header field to August 18, 2021 ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"February 11, 2021\"\n$find.Replacement.Text = \"August 18, 2021\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# --- 2. Insert a new \"South Korea\" hotline row right before the \"Taiwan\" row ---\n# (this is the APAC hotline table - the first table in the document)\n$t = $d.Tables(1)\n\n$taiwanRow = 0\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    if ($t.Cell($i, 1).Range.Text -like \"Taiwan*\") {\n        $taiwanRow = $i\n        break\n    }\n}\n\nif ($taiwanRow -gt 0) {\n    # Add the new row *after* Taiwan so it inherits Taiwan's (non-highlighted) row\n    # formatting instead of the preceding (highlighted) Singapore row, then swap\n    # the text so the new row ends up positioned before Taiwan.\n    $t.Rows.Add($t.Rows($taiwanRow + 1)) | Out-Null\n\n    $t.Cell($taiwanRow + 1, 1).Range.Text = $t.Cell($taiwanRow, 1).Range.Text\n    $t.Cell($taiwanRow + 1, 2).Range.Text = $t.Cell($taiwanRow, 2).Range.Text\n\n    $t.Cell($taiwanRow, 1).Range.Text = \"South Korea\"\n    $t.Cell($taiwanRow, 2).Range.Text = \"(+82) 798 611 4890\"\n}\n\n# --- 3. Simplify the \"Vietnam (Mobifone)\" country cell to a single run ---\n# (collapses the proofErr-wrapped, multi-run \"Vietnam (\" + \"Mobifone\" + \")\" text)\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Vietnam (Mobifone)\"\n$find2.Replacement.Text = \"Vietnam (Mobifone)\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
